$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.200971579659949
$ws.Range("C2").Value = 0.3361774439583201
$ws.Range("E2").Value = 0.4255972566567436
$ws.Range("F2").Value = 0.4443680307746121
$ws.Range("G2").Value = 0.002368582064315912
$ws.Range("O2").Value = 1.069890387247725
$ws.Range("B3").Value = 1.051621418813795
$ws.Range("C3").Value = 0.2967115533359959
$ws.Range("E3").Value = 0.3711417644439621
$ws.Range("F3").Value = 0.3878228170618172
$ws.Range("G3").Value = 0.002371829086850574
$ws.Range("O3").Value = 1.090993658087555
$ws.Range("B4").Value = 0.959670628208471
$ws.Range("C4").Value = 0.2723734315802346
$ws.Range("E4").Value = 0.3378068073238722
$ws.Range("F4").Value = 0.3531389305169483
$ws.Range("G4").Value = 0.002373923396827502
$ws.Range("O4").Value = 1.10583627643517
$ws.Range("B5").Value = 0.9221387855381522
$ws.Range("C5").Value = 0.2624293250775906
$ws.Range("E5").Value = 0.3242458192022468
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 0.002374802232908868
$ws.Range("O5").Value = 1.112355092281305
$ws.Range("B6").Value = 0.9159030037884008
$ws.Range("C6").Value = 0.2607765510839499
$ws.Range("E6").Value = 0.3219953722615969
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 0.002374949698750221
$ws.Range("O6").Value = 1.11346583505663
$ws.Range("B7").Value = 0.9591647062063657
$ws.Range("C7").Value = 0.2722394269730728
$ws.Range("E7").Value = 0.3376238275949675
$ws.Range("F7").Value = 0.3529483938344953
$ws.Range("G7").Value = 0.002373935146127418
$ws.Range("O7").Value = 1.105922292040518
$ws.Range("B8").Value = 1.149528212580606
$ws.Range("C8").Value = 0.3225918739981921
$ws.Range("E8").Value = 0.4067987675111766
$ws.Range("F8").Value = 0.4248636149813478
$ws.Range("G8").Value = 0.002369680804228347
$ws.Range("O8").Value = 1.076773310276039
$ws.Range("B9").Value = 1.520806752169165
$ws.Range("C9").Value = 0.420475497740938
$ws.Range("E9").Value = 0.5433530667459223
$ws.Range("F9").Value = 0.5661985755041457
$ws.Range("G9").Value = 0.002362132564603717
$ws.Range("O9").Value = 1.034736413520363
$ws.Range("B10").Value = 1.79231821255172
$ws.Range("C10").Value = 0.4918526480182663
$ws.Range("E10").Value = 0.6443799984424032
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 0.002357065726569811
$ws.Range("O10").Value = 1.013301976229172
$ws.Range("B11").Value = 1.915556186411379
$ws.Range("C11").Value = 0.5242045864798683
$ws.Range("E11").Value = 0.6905252142550324
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 0.00235486349148889
$ws.Range("O11").Value = 1.005652523673376
$ws.Range("B12").Value = 1.962182871533628
$ws.Range("C12").Value = 0.5364381167842112
$ws.Range("E12").Value = 0.7080287291752114
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 0.00235404424023058
$ws.Range("O12").Value = 1.003062099067137
$ws.Range("B13").Value = 1.952142820816562
$ws.Range("C13").Value = 0.5338041908304376
$ws.Range("E13").Value = 0.7042576920939752
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 0.002354220028606845
$ws.Range("O13").Value = 1.003606306865947
$ws.Range("B14").Value = 1.919393019016013
$ws.Range("C14").Value = 0.5252113988251494
$ws.Range("E14").Value = 0.6919646376862545
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 0.002354795797299052
$ws.Range("O14").Value = 1.00543324643553
$ws.Range("B15").Value = 1.899327455705759
$ws.Range("C15").Value = 0.5199457790796487
$ws.Range("E15").Value = 0.6844386760425749
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 0.002355150382586546
$ws.Range("O15").Value = 1.006592305637696
$ws.Range("B16").Value = 1.784258682090524
$ws.Range("C16").Value = 0.4897359579907743
$ws.Range("E16").Value = 0.6413682745653944
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 0.002357211707493911
$ws.Range("O16").Value = 1.013844518019482
$ws.Range("B17").Value = 1.71359649292765
$ws.Range("C17").Value = 0.4711726616582723
$ws.Range("E17").Value = 0.6149956799368539
$ws.Range("F17").Value = 0.6400460337215605
$ws.Range("G17").Value = 0.002358502509302529
$ws.Range("O17").Value = 1.018834597268636
$ws.Range("B18").Value = 1.672927769901605
$ws.Range("C18").Value = 0.4604844774078174
$ws.Range("E18").Value = 0.5998443722030942
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 0.002359254615137588
$ws.Range("O18").Value = 1.021902285317395
$ws.Range("B19").Value = 1.659153670387752
$ws.Range("C19").Value = 0.4568637563723996
$ws.Range("E19").Value = 0.5947173361642513
$ws.Range("F19").Value = 0.6191636801734006
$ws.Range("G19").Value = 0.002359510928752079
$ws.Range("O19").Value = 1.022974752150645
$ws.Range("B20").Value = 1.721121274278232
$ws.Range("C20").Value = 0.4731499052709864
$ws.Range("E20").Value = 0.6178012583380905
$ws.Range("F20").Value = 0.642933953830422
$ws.Range("G20").Value = 0.002358364100873137
$ws.Range("O20").Value = 1.018282921909332
$ws.Range("B21").Value = 1.92901355009576
$ws.Range("C21").Value = 0.5277357875060602
$ws.Range("E21").Value = 0.6955745910188114
$ws.Range("F21").Value = 0.7228739723492197
$ws.Range("G21").Value = 0.002354626282061364
$ws.Range("O21").Value = 1.004888285432571
$ws.Range("B22").Value = 2.064644352272865
$ws.Range("C22").Value = 0.5633089043030282
$ws.Range("E22").Value = 0.7465759933968314
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 0.002352268973703608
$ws.Range("O22").Value = 0.997921135936565
$ws.Range("B23").Value = 1.992277965507355
$ws.Range("C23").Value = 0.5443323474696058
$ws.Range("E23").Value = 0.7193390579660104
$ws.Range("F23").Value = 0.7472568307830727
$ws.Range("G23").Value = 0.002353519310095399
$ws.Range("O23").Value = 1.001474760067936
$ws.Range("B24").Value = 1.717719457852638
$ws.Range("C24").Value = 0.472256042919355
$ws.Range("E24").Value = 0.6165328233760761
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 0.002358426644270412
$ws.Range("O24").Value = 1.0185317153298
$ws.Range("B25").Value = 1.420586322052145
$ws.Range("C25").Value = 0.3940888007044236
$ws.Range("E25").Value = 0.5062994664214244
$ws.Range("F25").Value = 0.5279251897347166
$ws.Range("G25").Value = 0.002364090075321121
$ws.Range("O25").Value = 1.044465640841835
